$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: "SyncJobDAOImpl - create (), update (), find ()"
#           -> "SyncJobDAOImpl -findAll ()"  (split across 4 runs)
# ------------------------------------------------------------------
$dash = [char]8211

$daoPara = $null
foreach ($pp in $d.Paragraphs) {
    if ($pp.Range.Text -match "SyncJobDAOImpl") {
        $daoPara = $pp
    }
}

$target = $d.Content
$target.Start = $daoPara.Range.Start
$target.End = $daoPara.Range.End
$target.Find.Execute("SyncJobDAOImpl " + $dash + " create (), update (), find ()")

$rng = $d.Range($target.Start, $target.End)
$rng.Text = ""
$pos = $rng.Start

$r1 = $d.Range($pos, $pos)
$r1.InsertAfter("SyncJobDAOImpl " + $dash)
$r1.Font.Size = 12
$pos = $r1.End

$r2 = $d.Range($pos, $pos)
$r2.InsertAfter("find")
$r2.Font.Size = 12
$pos = $r2.End

$r3 = $d.Range($pos, $pos)
$r3.InsertAfter("All")
$r3.Font.Size = 12
$pos = $r3.End

$r4 = $d.Range($pos, $pos)
$r4.InsertAfter(" ()")
$r4.Font.Size = 12
$pos = $r4.End

# ------------------------------------------------------------------
# Change 2: insert a highlighted "Note: ..." paragraph (plus a
#           trailing blank paragraph) right after the blank line
#           that follows the Sync-Job-Status table, before
#           "Class Diagram:".
# ------------------------------------------------------------------
$blankPara = $null
foreach ($pp in $d.Paragraphs) {
    if (($pp.Range.Text -eq [char]13) -and ($pp.Next -ne $null) -and ($pp.Next.Range.Text -match "^Class Diagram")) {
        $blankPara = $pp
    }
}

$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$quoteOpen = [char]8220
$quoteClose = [char]8221
$noteText = "Note: " + $quoteOpen + "validateSyncStatus()" + $quoteClose + " method from Base controller is pluggable from any other controller as required."

$xml = "<w:p $ns w:rsidR='00B610D6' w:rsidRDefault='00B610D6' w:rsidP='00B610D6'><w:pPr><w:rPr><w:sz w:val='28'/><w:szCs w:val='28'/></w:rPr></w:pPr></w:p>" + `
       "<w:p $ns><w:pPr><w:rPr><w:sz w:val='24'/><w:szCs w:val='28'/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val='24'/><w:szCs w:val='28'/><w:highlight w:val='yellow'/></w:rPr><w:t>$noteText</w:t></w:r></w:p>" + `
       "<w:p $ns><w:pPr><w:rPr><w:sz w:val='24'/><w:szCs w:val='28'/></w:rPr></w:pPr></w:p>"

$blankPara.Range.InsertXML($xml)

# ------------------------------------------------------------------
# Change 3: remove the stray "      " (whitespace-only) paragraph
#           that sits between the Sequence-Diagram hyperlink
#           paragraph and the _GoBack bookmark paragraph.
# ------------------------------------------------------------------
$spacesPara = $null
foreach ($pp in $d.Paragraphs) {
    $t = $pp.Range.Text
    if ($t -match "^[ ]{3,}\r$") {
        $spacesPara = $pp
    }
}
if ($spacesPara -ne $null) {
    $spacesPara.Range.Delete()
}
